# Rename the inline logo pictures embedded in the document's headers/footers.
#
# The document has a title page (titlePg), so the section exposes distinct
# "primary" and "first page" header/footer stories:
#   - Footer(wdHeaderFooterPrimary)   -> Pearson logo, currently "image1.png"
#   - Footer(wdHeaderFooterFirstPage) -> Pearson logo, currently "image1.png"
#   - Header(wdHeaderFooterFirstPage) -> BTec logo,    currently "image2.jpg"
#
# We need to swap the display names:
#   Pearson logos: image1.png -> image2.png
#   BTec logo:     image2.jpg -> image1.jpg
#
# InlineShape has no settable Name property (matches real Word), so each
# picture is promoted to a floating Shape (which does expose Name), renamed,
# then converted back to an inline picture in place.

$wdHeaderFooterPrimary = 1
$wdHeaderFooterFirstPage = 2

function Rename-LogoInlineShape($story, $newName) {
    $inline = $story.Range.InlineShapes.Item(1)
    $shape = $inline.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

$d = $word.ActiveDocument
$section = $d.Sections(1)

# Primary footer - Pearson logo
Rename-LogoInlineShape $section.Footers($wdHeaderFooterPrimary) "image2.png"

# First-page footer - Pearson logo
Rename-LogoInlineShape $section.Footers($wdHeaderFooterFirstPage) "image2.png"

# First-page header - BTec logo
Rename-LogoInlineShape $section.Headers($wdHeaderFooterFirstPage) "image1.jpg"
